$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume/name/link data per the Sep 23 2023 refresh.
# Some Price (column D) values look like plain numbers (e.g. '211.49'); Excel's
# own type-inference would silently coerce those to numeric cells (dropping
# trailing zeros, changing storage type) unlike the source data, which keeps them
# as text. Force text storage for those specific cells, then strip the temporary
# number-format override so no stray style is left behind on the cell.
$ws.Range("D2").Value = '26.719.53'
$ws.Range("D3").Value = '1.599.14'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.49'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("E6").Value = '  -0.61%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0619'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").Value = '1.824.01'
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.631.94'
$ws.Range("E13").Value = '  +2.13%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.04'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.33'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '26.688.62'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("E18").Value = '  +4.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '209.88'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("E21").Value = '  +3.83%  '
$ws.Range("E22").Value = '  +0.78%  '
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.93'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.31'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("E30").Value = '  +3.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.25'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("E33").Value = '  +1.79%  '
$ws.Range("D34").Value = '1.289.46'
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("E35").Value = '  -5.41%  '
$ws.Range("E36").Value = '  +0.93%  '
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("E39").Value = '  +17.50%  '
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.785'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.13%  '
$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.10'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.78%  '
$ws.Range("D45").Value = '1.736.76'
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.36'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.89%  '
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.100'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0508'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.65%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.37'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.94%  '
